$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# Replace "Stats" with "News" in the ViewModels column for the "Top bar" widget sections
$ws.Range("D14").Value = "News"
$ws.Range("D15").Value = "News"
$ws.Range("D16").Value = "News"
$ws.Range("D19").Value = "News"
$ws.Range("D20").Value = "News"
$ws.Range("D21").Value = "News"

# Move the active cell selection (cosmetic)
$ws.Range("H13").Select()
